# ---------------------------------------------------------------------------
# "added the MSR site summary"
#
# The edit:
#   1. Drops the stray "_GoBack" bookmark that used to sit in the middle of
#      the "...understanding the Rails structure..." sentence.
#   2. Appends a new "Master of Science in Robotics Site" section (a bold
#      heading followed by two body paragraphs, separated by blank
#      paragraphs) to the end of the document. The "_GoBack" bookmark
#      reappears inside the second new paragraph, right after
#      the Jekyll “back-end” (not really) aside -- i.e. it moved rather
#      than simply being deleted.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# 1. Remove the old "_GoBack" bookmark wherever Word last left it.
try {
    $d.Bookmarks("_GoBack").Delete()
} catch {
    # Nothing to remove -- fine, the target state has no bookmark here either.
}

# 2. Insert the new section (7 paragraphs) right after the final paragraph
#    ("...using The Book."), before the sectPr. Using Range.InsertXML lets us
#    place exact OOXML (including the empty <w:p/> spacer paragraphs and the
#    relocated "_GoBack" bookmark) instead of reconstructing it run-by-run
#    through the high-level Range/Selection typing API.
$tail = $d.Content
$tail.Collapse(0)  # wdCollapseEnd -- move to the very end of the body content
$newSectionXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/><w:p/><w:p><w:r><w:rPr><w:b/></w:rPr><w:t>Master of Science in Robotics Site</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">After working for a summer to develop projects and curriculum for one of the courses in the new Master of Science in Robotics Program at Northwestern University, </w:t></w:r><w:r><w:t>my boss asked me to build the new</w:t></w:r><w:r><w:t xml:space="preserve"> program’s </w:t></w:r><w:r><w:t>website.</w:t></w:r><w:r><w:t xml:space="preserve"> His vision was to have </w:t></w:r><w:r><w:t xml:space="preserve">a site that </w:t></w:r><w:r><w:t>adequately</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">displayed </w:t></w:r><w:r><w:t xml:space="preserve">everything that the program’s students </w:t></w:r><w:r><w:t>would be working on</w:t></w:r><w:r><w:t xml:space="preserve"> in addition to housing </w:t></w:r><w:r><w:t xml:space="preserve">a database of </w:t></w:r><w:r><w:t>resources</w:t></w:r><w:r><w:t xml:space="preserve"> that students and faculty would contribute to over time.</w:t></w:r><w:r><w:t xml:space="preserve"> In addition, he wanted</w:t></w:r><w:r><w:t xml:space="preserve"> a website that </w:t></w:r><w:r><w:t xml:space="preserve">could be easily maintained </w:t></w:r><w:r><w:t>using GitH</w:t></w:r><w:r><w:t>ub</w:t></w:r><w:r><w:t xml:space="preserve"> and M</w:t></w:r><w:r><w:t>arkdown.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">To </w:t></w:r><w:r><w:t xml:space="preserve">meet </w:t></w:r><w:r><w:t xml:space="preserve">that last requirement in particular, I chose to </w:t></w:r><w:r><w:t xml:space="preserve">develop the website using </w:t></w:r><w:r><w:t xml:space="preserve">Jekyll, a </w:t></w:r><w:r><w:t>“</w:t></w:r><w:r><w:t>back-end</w:t></w:r><w:r><w:t>” (not really)</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>static-site generator.</w:t></w:r><w:r><w:t xml:space="preserve"> Jekyll’s main strengths</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>include</w:t></w:r><w:r><w:t xml:space="preserve"> its</w:t></w:r><w:r><w:t xml:space="preserve"> incredibly </w:t></w:r><w:r><w:t>painless</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">integration </w:t></w:r><w:r><w:t>in</w:t></w:r><w:r><w:t>to</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">the </w:t></w:r><w:r><w:t>GitHub Pages</w:t></w:r><w:r><w:t xml:space="preserve"> hosting</w:t></w:r><w:r><w:t xml:space="preserve"> service </w:t></w:r><w:r><w:t xml:space="preserve">and </w:t></w:r><w:r><w:t>its</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>feature that allows a site’s content to be maintained in Markdown.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$tail.InsertXML($newSectionXml)
